$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product name in A3 from "Smart TV" to "Mobiles 5g"
$ws.Range("A3").Value = "Mobiles 5g"

# Move the active cell selection from A7 to A6
$ws.Range("A6").Select()
